# Auto-generated PowerShell/Excel COM-interop script
# Applies numeric corrections to Leviathan_Profits sheets as described in the commit diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 566.5
$ws.Range("I18").Value = 566.5
$ws.Range("K18").Value = 566.5
$ws.Range("M18").Value = -282.5
$ws.Range("H33").Value = 100.0
$ws.Range("I33").Value = 0.0
$ws.Range("K33").Value = 0.0
$ws.Range("M33").ClearContents()
$ws.Range("H86").Value = 2600.1538
$ws.Range("I86").Value = 2580.2
$ws.Range("J86").Value = 2666.6667
$ws.Range("K86").Value = 2580.2
$ws.Range("L86").Value = 2666.6667
$ws.Range("M86").Value = -1457.2
$ws.Range("N86").Value = -4912.6667
$ws.Range("H89").Value = 2600.1538
$ws.Range("I89").Value = 2580.2
$ws.Range("J89").Value = 2666.6667
$ws.Range("K89").Value = 12901.0
$ws.Range("L89").Value = 13333.3335
$ws.Range("M89").Value = -7285.0
$ws.Range("N89").Value = -24565.3335
$ws.Range("H111").Value = 996.8571
$ws.Range("I111").Value = 996.5
$ws.Range("J111").Value = 999.0
$ws.Range("K111").Value = 2989.5
$ws.Range("L111").Value = 2997.0
$ws.Range("M111").Value = 77.5
$ws.Range("N111").Value = -9131.0
$ws.Range("H113").Value = 4465.9375
$ws.Range("J113").Value = 4830.364
$ws.Range("L113").Value = 4830.364
$ws.Range("N113").Value = -11338.364
$ws.Range("H116").Value = 4999.8
$ws.Range("I116").Value = 4999.0
$ws.Range("K116").Value = 4999.0
$ws.Range("M116").Value = -1557.0
$ws.Range("H132").Value = 3747.111
$ws.Range("I132").Value = 1496.0
$ws.Range("K132").Value = 4488.0
$ws.Range("M132").Value = -1958.0
$ws.Range("H137").Value = 2995.8462
$ws.Range("I137").Value = 1749.5555
$ws.Range("J137").Value = 5800.0
$ws.Range("K137").Value = 5248.666499999999
$ws.Range("L137").Value = 17400.0
$ws.Range("M137").Value = -2698.666499999999
$ws.Range("N137").Value = -22500.0
$ws.Range("H138").Value = 3336.72
$ws.Range("J138").Value = 3541.366
$ws.Range("L138").Value = 10624.098
$ws.Range("N138").Value = -20904.098
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1368.5217
$ws.Range("I2").Value = 1294.1111
$ws.Range("K2").Value = 1294.1111
$ws.Range("M2").Value = -1181.1111
$ws.Range("H32").Value = 5569.683
$ws.Range("I32").Value = 3967.029
$ws.Range("J32").Value = 14076.077
$ws.Range("K32").Value = 3967.029
$ws.Range("L32").Value = 14076.077
$ws.Range("M32").Value = -3680.029
$ws.Range("N32").Value = -14650.077
$ws.Range("H45").Value = 7306.2383
$ws.Range("I45").Value = 10092.833
$ws.Range("J45").Value = 3590.7778
$ws.Range("K45").Value = 10092.833
$ws.Range("L45").Value = 3590.7778
$ws.Range("M45").Value = -9715.833
$ws.Range("N45").Value = -4344.7778
$ws.Range("H116").Value = 1368.5217
$ws.Range("I116").Value = 1294.1111
$ws.Range("K116").Value = 1294.1111
$ws.Range("M116").Value = 999.8888999999999
$ws.Range("H124").Value = 68073.0
$ws.Range("J124").Value = 68073.0
$ws.Range("L124").Value = 68073.0
$ws.Range("N124").Value = -77893.0
$ws.Range("H132").Value = 2396.6086
$ws.Range("I132").Value = 2006.6923
$ws.Range("J132").Value = 4569.0
$ws.Range("K132").Value = 6020.0769
$ws.Range("L132").Value = 13707.0
$ws.Range("M132").Value = -3490.0769
$ws.Range("N132").Value = -18767.0
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1368.5217
$ws.Range("I3").Value = 1294.1111
$ws.Range("K3").Value = 1294.1111
$ws.Range("M3").Value = -1180.1111
$ws.Range("H22").Value = 1552.0
$ws.Range("I22").Value = 1552.0
$ws.Range("K22").Value = 1552.0
$ws.Range("M22").Value = -1379.0
$ws.Range("H94").Value = 32259444.0
$ws.Range("H99").Value = 2069.3333
$ws.Range("I99").Value = 2085.9167
$ws.Range("K99").Value = 2085.9167
$ws.Range("M99").Value = -587.9167000000002
$ws.Range("H107").Value = 2703.6667
$ws.Range("I107").Value = 2623.238
$ws.Range("K107").Value = 2623.238
$ws.Range("M107").Value = -703.2379999999998
$ws.Range("H117").Value = 40102.168
$ws.Range("J117").Value = 40102.168
$ws.Range("L117").Value = 40102.168
$ws.Range("N117").Value = -49280.168
$ws.Range("H134").Value = 216187.5
$ws.Range("I134").Value = 301861.2
$ws.Range("K134").Value = 905583.6000000001
$ws.Range("M134").Value = -903048.6000000001
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 100000320.0
$ws.Range("I7").Value = 200000350.0
$ws.Range("K7").Value = 200000350.0
$ws.Range("M7").Value = -200000237.0
$ws.Range("H15").Value = 219.2
$ws.Range("J15").Value = 295.0
$ws.Range("L15").Value = 295.0
$ws.Range("N15").Value = -635.0
$ws.Range("H19").Value = 950.1818
$ws.Range("I19").Value = 595.2
$ws.Range("K19").Value = 595.2
$ws.Range("M19").Value = -425.2
$ws.Range("H22").Value = 1000.0
$ws.Range("I22").Value = 1000.0
$ws.Range("K22").Value = 1000.0
$ws.Range("M22").Value = -650.0
$ws.Range("H24").Value = 950.1818
$ws.Range("I24").Value = 595.2
$ws.Range("K24").Value = 595.2
$ws.Range("M24").Value = -425.2
$ws.Range("H74").Value = 57807.5
$ws.Range("J74").Value = 57807.5
$ws.Range("L74").Value = 57807.5
$ws.Range("N74").Value = -59555.5
$ws.Range("H77").Value = 57807.5
$ws.Range("J77").Value = 57807.5
$ws.Range("L77").Value = 173422.5
$ws.Range("N77").Value = -182158.5
$ws.Range("H107").Value = 1948.4584
$ws.Range("I107").Value = 1413.6
$ws.Range("K107").Value = 1413.6
$ws.Range("M107").Value = 506.4000000000001
$ws.Range("H132").Value = 3796.647
$ws.Range("I132").Value = 3741.8462
$ws.Range("K132").Value = 11225.5386
$ws.Range("M132").Value = -8695.5386
$ws.Range("H134").Value = 2400.5637
$ws.Range("I134").Value = 2025.1
$ws.Range("K134").Value = 6075.299999999999
$ws.Range("M134").Value = -3540.299999999999
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 41.8
$ws.Range("I7").Value = 41.8
$ws.Range("K7").Value = 125.4
$ws.Range("M7").Value = -13.39999999999999
$ws.Range("H114").Value = 18183378.0
$ws.Range("J114").Value = 1837.5
$ws.Range("L114").Value = 5512.5
$ws.Range("N114").Value = -12020.5
$ws.Range("H121").Value = 19668840.0
$ws.Range("I121").Value = 66666864.0
$ws.Range("J121").Value = 86331.164
$ws.Range("K121").Value = 200000592.0
$ws.Range("L121").Value = 258993.492
$ws.Range("M121").Value = -199999282.0
$ws.Range("N121").Value = -261613.492
$ws.Range("H129").Value = 76756.11
$ws.Range("I129").Value = 401064.2
$ws.Range("J129").Value = 3049.7273
$ws.Range("K129").Value = 1203192.6
$ws.Range("L129").Value = 9149.1819
$ws.Range("M129").Value = -1198192.6
$ws.Range("N129").Value = -19149.1819
$ws.Range("H139").Value = 3481.1428
$ws.Range("I139").Value = 3481.1428
$ws.Range("K139").Value = 10443.4284
$ws.Range("M139").Value = -5303.428400000001
$ws.Range("H141").Value = 3119.8667
$ws.Range("I141").Value = 2771.2856
$ws.Range("J141").Value = 8000.0
$ws.Range("K141").Value = 8313.856800000001
$ws.Range("L141").Value = 24000.0
$ws.Range("M141").Value = -3133.856800000001
$ws.Range("N141").Value = -34360.0
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 4351.457
$ws.Range("I132").Value = 3146.44
$ws.Range("J132").Value = 7364.0
$ws.Range("K132").Value = 9439.32
$ws.Range("L132").Value = 22092.0
$ws.Range("M132").Value = -6909.32
$ws.Range("N132").Value = -27152.0
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 3616.5173
$ws.Range("I22").Value = 3990.7917
$ws.Range("K22").Value = 3990.7917
$ws.Range("M22").Value = -3695.7917
$ws.Range("H27").Value = 3616.5173
$ws.Range("I27").Value = 3990.7917
$ws.Range("K27").Value = 3990.7917
$ws.Range("M27").Value = -3883.7917
$ws.Range("H46").Value = 3592.9546
$ws.Range("I46").Value = 3795.8333
$ws.Range("K46").Value = 3795.8333
$ws.Range("M46").Value = -3607.8333
$ws.Range("H100").Value = 4926.154
$ws.Range("I100").Value = 3459.0667
$ws.Range("K100").Value = 3459.0667
$ws.Range("M100").Value = -2918.0667
$ws.Range("H136").Value = 6081.85
$ws.Range("I136").Value = 5878.793
$ws.Range("J136").Value = 6617.1816
$ws.Range("K136").Value = 17636.379
$ws.Range("L136").Value = 19851.5448
$ws.Range("M136").Value = -15086.379
$ws.Range("N136").Value = -24951.5448
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H52").Value = 25498.334
$ws.Range("J52").Value = 25498.334
$ws.Range("L52").Value = 25498.334
$ws.Range("N52").Value = -25950.334
$ws.Range("H135").Value = 48738.0
$ws.Range("J135").Value = 48738.0
$ws.Range("L135").Value = 48738.0
$ws.Range("N135").Value = -58878.0
